$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9773720502853394
$ws.Range("B1").Value = 1.877420425415039
$ws.Range("C1").Value = 6.192467212677002
$ws.Range("D1").Value = 3.599578619003296
$ws.Range("E1").Value = 1.325616359710693
